# Actualización desde MV -datos-
# Adds two new daily rows (05-11-2021 and 08-11-2021) to the monetary policy
# rate table, and fills in the previously-missing C/D values for the two
# preceding rows (03-11-2021 and 04-11-2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete the existing rows that were missing the facilidad permanente
# (liquidity / deposit) columns C and D.
$ws.Range("C211").Value = 3
$ws.Range("D211").Value = 2.5

$ws.Range("C212").Value = 3
$ws.Range("D212").Value = 2.5

# New row for 05-11-2021
# (NumberFormat "@" forces the date-looking string to be stored as literal
# text instead of being auto-converted to a date serial; Style is reset back
# to Normal afterwards so the cell keeps the workbook's default formatting.)
$ws.Range("A213").NumberFormat = "@"
$ws.Range("A213").Value = "05-11-2021"
$ws.Range("A213").Style = "Normal"
$ws.Range("B213").Value = 2.75
$ws.Range("C213").Value = 3
$ws.Range("D213").Value = 2.5

# New row for 08-11-2021 (no TPM value published for this date)
$ws.Range("A214").NumberFormat = "@"
$ws.Range("A214").Value = "08-11-2021"
$ws.Range("A214").Style = "Normal"
$ws.Range("C214").Value = 3
$ws.Range("D214").Value = 2.5
